# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Modesto" variety at row 48,
# pushing the existing data down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 48.
$ws.Rows("48:50").Insert()

# Shared boilerplate values (identical across every data row in this sheet).
$marketId = 8
$market   = "Terminal La Palmera de La Serena"
$region   = "Coquimbo"
$codreg   = 4
$tipo     = "Fruta"
$prodId   = 100103
$prod     = "Frutos de hueso (carozo)"
$catId    = 100103003
$cat      = "Damasco"

function Set-DataRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidad
    )

    $ws.Cells.Item($Row, 1).Value  = $marketId
    $ws.Cells.Item($Row, 2).Value  = $market
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $prodId
    $ws.Cells.Item($Row, 8).Value  = $prod
    $ws.Cells.Item($Row, 9).Value  = $catId
    $ws.Cells.Item($Row, 10).Value = $cat
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

Set-DataRow 48 44944 "Modesto" "Especial" 200 24000 25000 24500 "$/caja 16 kilos" "Región de O'Higgins" 1531 16
Set-DataRow 49 44944 "Modesto" "Primera"  240 21000 22000 21500 "$/caja 16 kilos" "Región de O'Higgins" 1344 16
Set-DataRow 50 44944 "Modesto" "Segunda"  200 18000 19000 18500 "$/caja 16 kilos" "Región de O'Higgins" 1156 16
